$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "L1cam"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 22.59487733333333
$ws.Range("H2").Value = 67.784632
$ws.Range("I2").Value = 0.7395019553569895
$ws.Range("J2").Value = 0.7395019553569895
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.534538333333333
$ws.Range("N2").Value = 4.603615
$ws.Range("O2").Value = 0.1494637976135089
$ws.Range("P2").Value = 0.1494637976135089
$ws.Range("Q2").Value = 34.67270540496444
$ws.Range("R2").Value = 312.05434864468
$ws.Range("S2").Value = 0.1105287705902712
$ws.Range("T2").Value = 0.1105287705902712

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "L1cam"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 22.59487733333333
$ws.Range("H3").Value = 67.784632
$ws.Range("I3").Value = 0.7395019553569895
$ws.Range("J3").Value = 0.7395019553569895
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.973328333333334
$ws.Range("N3").Value = 14.919985
$ws.Range("O3").Value = 0.4844014146353658
$ws.Range("P3").Value = 0.4844014146353658
$ws.Range("Q3").Value = 112.3717436300578
$ws.Range("R3").Value = 1011.34569267052
$ws.Range("S3").Value = 0.3582157933005449
$ws.Range("T3").Value = 0.3582157933005448

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "L1cam"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 22.59487733333333
$ws.Range("H4").Value = 67.784632
$ws.Range("I4").Value = 0.7395019553569895
$ws.Range("J4").Value = 0.7395019553569895
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.75909
$ws.Range("N4").Value = 11.27727
$ws.Range("O4").Value = 0.3661347877511252
$ws.Range("P4").Value = 0.3661347877511252
$ws.Range("Q4").Value = 84.93617743496
$ws.Range("R4").Value = 764.42559691464
$ws.Range("S4").Value = 0.2707573914661734
$ws.Range("T4").Value = 0.2707573914661734

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "L1cam"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3045986666666667
$ws.Range("H5").Value = 0.913796
$ws.Range("I5").Value = 0.00996913177602551
$ws.Range("J5").Value = 0.00996913177602551
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.534538333333333
$ws.Range("N5").Value = 4.603615
$ws.Range("O5").Value = 0.1494637976135089
$ws.Range("P5").Value = 0.1494637976135089
$ws.Range("Q5").Value = 0.4674183302822222
$ws.Range("R5").Value = 4.20676497254
$ws.Range("S5").Value = 0.001490024294154278
$ws.Range("T5").Value = 0.001490024294154278

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "L1cam"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3045986666666667
$ws.Range("H6").Value = 0.913796
$ws.Range("I6").Value = 0.00996913177602551
$ws.Range("J6").Value = 0.00996913177602551
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.973328333333334
$ws.Range("N6").Value = 14.919985
$ws.Range("O6").Value = 0.4844014146353658
$ws.Range("P6").Value = 0.4844014146353658
$ws.Range("Q6").Value = 1.514869179228889
$ws.Range("R6").Value = 13.63382261306
$ws.Range("S6").Value = 0.004829061534993134
$ws.Range("T6").Value = 0.004829061534993134

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "L1cam"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3045986666666667
$ws.Range("H7").Value = 0.913796
$ws.Range("I7").Value = 0.00996913177602551
$ws.Range("J7").Value = 0.00996913177602551
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.75909
$ws.Range("N7").Value = 11.27727
$ws.Range("O7").Value = 0.3661347877511252
$ws.Range("P7").Value = 0.3661347877511252
$ws.Range("Q7").Value = 1.14501380188
$ws.Range("R7").Value = 10.30512421692
$ws.Range("S7").Value = 0.003650045946878098
$ws.Range("T7").Value = 0.003650045946878098

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "L1cam"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.654706
$ws.Range("H8").Value = 22.964118
$ws.Range("I8").Value = 0.2505289128669849
$ws.Range("J8").Value = 0.2505289128669849
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.534538333333333
$ws.Range("N8").Value = 4.603615
$ws.Range("O8").Value = 0.1494637976135089
$ws.Range("P8").Value = 0.1494637976135089
$ws.Range("Q8").Value = 11.74643978739667
$ws.Range("R8").Value = 105.71795808657
$ws.Range("S8").Value = 0.03744500272908345
$ws.Range("T8").Value = 0.03744500272908345

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "L1cam"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.654706
$ws.Range("H9").Value = 22.964118
$ws.Range("I9").Value = 0.2505289128669849
$ws.Range("J9").Value = 0.2505289128669849
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.973328333333334
$ws.Range("N9").Value = 14.919985
$ws.Range("O9").Value = 0.4844014146353658
$ws.Range("P9").Value = 0.4844014146353658
$ws.Range("Q9").Value = 38.06936623313667
$ws.Range("R9").Value = 342.62429609823
$ws.Range("S9").Value = 0.1213565597998278
$ws.Range("T9").Value = 0.1213565597998278

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "L1cam"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.654706
$ws.Range("H10").Value = 22.964118
$ws.Range("I10").Value = 0.2505289128669849
$ws.Range("J10").Value = 0.2505289128669849
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.75909
$ws.Range("N10").Value = 11.27727
$ws.Range("O10").Value = 0.3661347877511252
$ws.Range("P10").Value = 0.3661347877511252
$ws.Range("Q10").Value = 28.77472877754
$ws.Range("R10").Value = 258.97255899786
$ws.Range("S10").Value = 0.09172735033807365
$ws.Range("T10").Value = 0.09172735033807365
